# Updated transaction file in 27.04.2017
#
# poloniex history.xlsx - "Sheet1" transaction log. Two orders that were
# previously "IN PROGRESS" (ZEC buy & NXT buy, rows 24/25) got cancelled and
# two new orders for the same pairs were placed; a pending XRP sell (row 27)
# completed. Two brand-new rows (28/29) record the replacement ZEC/NXT buys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GREEN = 5287936   # RGB(0,176,80) packed BGR, used for "Buy"/"~x%" runs

function Set-RowHeight($row, $height) {
    $ws.Rows.Item($row).RowHeight = $height
}

# ---------------------------------------------------------------------------
# Row 24: ZEC buy order -> CANCEL, stamped with a finalize date/time
# ---------------------------------------------------------------------------
$ws.Range("H24").Value = "CANCEL"
$ws.Range("I24").Value = "2017-04-27 09:40:37  `n"
Set-RowHeight 24 14.25

# ---------------------------------------------------------------------------
# Row 25: NXT buy order -> CANCEL, stamped with a finalize date/time
# ---------------------------------------------------------------------------
$ws.Range("H25").Value = "CANCEL"
$ws.Range("I25").Value = "4/27/2017  9:40:37 AM`n"
Set-RowHeight 25 14.25

# ---------------------------------------------------------------------------
# Row 27: XRP sell order -> DONE, with fee/profit/duration info filled in
# ---------------------------------------------------------------------------
$ws.Range("H27").Value = "DONE"
$ws.Range("I27").Value = "2017-04-26 18:59:06`n"
$ws.Range("J27").Value = "0.00974957 USDT (0.15%)"
$k27 = $ws.Range("K27")
$k27.Value = "     ~1.5%"
$k27.Characters(5, 6).Font.Color = $GREEN
$ws.Range("L27").Value = " 1 day"
Set-RowHeight 27 14.25

# ---------------------------------------------------------------------------
# Row 28 (new): another ZEC buy order, still in progress
# ---------------------------------------------------------------------------
$ws.Range("A28").Value2 = 42852.40320601852
$ws.Range("A28").NumberFormat = "m/d/yy h:mm"

$b28 = $ws.Range("B28")
$b28.Value = "            Buy"
$b28.Characters(13, 3).Font.Color = $GREEN

$ws.Range("C28").Value = "        ZEC"
$ws.Range("D28").Value = 78.77000058
$ws.Range("E28").Value = "            77.5USDT"
$ws.Range("F28").Value = "        0.1 ZEC"
$ws.Range("G28").Value = " ZEC/USDT0000003"
$ws.Range("H28").Value = "IN PROGRESS"

$i28 = $ws.Range("I28")
$i28.Value = " `n"
$i28.NumberFormat = "m/d/yy h:mm"
$i28.WrapText = $true

Set-RowHeight 28 14.25

# ---------------------------------------------------------------------------
# Row 29 (new): another NXT buy order, still in progress
# ---------------------------------------------------------------------------
$a29 = $ws.Range("A29")
$a29.Value2 = 42852.40320601852
$a29.NumberFormat = "m/d/yy h:mm"
$a29.WrapText = $true

$b29 = $ws.Range("B29")
$b29.Value = "            Buy"
$b29.Characters(13, 3).Font.Color = $GREEN

$ws.Range("C29").Value = "        NXT"
$ws.Range("D29").Value = 0.02896926
$ws.Range("E29").Value = "           0.0285USDT"
$ws.Range("F29").Value = "          435 NXT"
$ws.Range("G29").Value = " NXT/USDT0000001"
$ws.Range("H29").Value = "IN PROGRESS"

$i29 = $ws.Range("I29")
$i29.NumberFormat = "m/d/yy h:mm"
$i29.WrapText = $true

Set-RowHeight 29 14.25

# ---------------------------------------------------------------------------
# View state: scroll so row 11 is at the top, and leave G32 selected
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A11"), $true)
$ws.Range("G32").Select() | Out-Null
